$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(4, 5, 7, 8, 10, 11, 13, 16, 17, 20, 26, 40, 41)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = "lipid/free"
}
